# Apply updated crypto price/volume data as per the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''61.774.21'
$ws.Range("E2").Value = '  -0.26%  '
$ws.Range("D3").Value = '''3.408.78'
$ws.Range("E3").Value = '  +0.03%  '
$ws.Range("E4").Value = '  -0.26%  '
$ws.Range("D5").Value = '''413.30'
$ws.Range("E5").Value = '  +1.21%  '
$ws.Range("D6").Value = '''129.33'
$ws.Range("E6").Value = '  +0.65%  '
$ws.Range("D7").Value = '''0.621'
$ws.Range("E7").Value = '  -2.90%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("E9").Value = '  -1.03%  '
$ws.Range("E10").Value = '  -4.50%  '
$ws.Range("D11").Value = '''42.66'
$ws.Range("E11").Value = '  +0.68%  '
$ws.Range("D12").Value = '''0.0000217'
$ws.Range("E12").Value = '  -0.93%  '
$ws.Range("E13").Value = '  +2.23%  '
$ws.Range("D14").Value = '''3.949.24'
$ws.Range("E14").Value = '  -0.22%  '
$ws.Range("D16").Value = '''20.43'
$ws.Range("E16").Value = '  -1.63%  '
$ws.Range("D17").Value = '''3.437.07'
$ws.Range("D18").Value = '''12.65'
$ws.Range("E18").Value = '  +4.04%  '
$ws.Range("E19").Value = '  +1.08%  '
$ws.Range("D20").Value = '''61.830.88'
$ws.Range("E20").Value = '  -0.10%  '
$ws.Range("D21").Value = '''481.75'
$ws.Range("E21").Value = '  +11.71%  '
$ws.Range("D22").Value = '''90.67'
$ws.Range("E22").Value = '  +0.93%  '
$ws.Range("E23").Value = '  +3.76%  '
$ws.Range("D24").Value = '''13.08'
$ws.Range("E24").Value = '  +0.48%  '
$ws.Range("E25").Value = '  +2.66%  '
$ws.Range("D26").Value = '''9.80'
$ws.Range("E26").Value = '  +11.25%  '
$ws.Range("D27").Value = '''33.03'
$ws.Range("E27").Value = '  -1.55%  '
$ws.Range("D28").Value = '''4.75'
$ws.Range("E28").Value = '  -0.44%  '
$ws.Range("E29").Value = '  +2.28%  '
$ws.Range("B30").Value = 'Cosmos'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D30").Value = '''11.86'
$ws.Range("E30").Value = '  -0.53%  '
$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").Value = '''2.65'
$ws.Range("E31").Value = '  -1.32%  '
$ws.Range("D32").Value = '''0.167'
$ws.Range("E32").Value = '  -1.48%  '
$ws.Range("E33").Value = '  -3.13%  '
$ws.Range("D34").Value = '''40.95'
$ws.Range("E34").Value = '  -3.37%  '
$ws.Range("D36").Value = '''58.26'
$ws.Range("E36").Value = '  +8.05%  '
$ws.Range("D37").Value = '''0.0486'
$ws.Range("E37").Value = '  -2.23%  '
$ws.Range("E38").Value = '  +0.06%  '
$ws.Range("D39").Value = '''3.03'
$ws.Range("E39").Value = '  +4.42%  '
$ws.Range("D40").Value = '''149.10'
$ws.Range("E40").Value = '  +5.95%  '
$ws.Range("D41").Value = '''0.323'
$ws.Range("E41").Value = '  +3.64%  '
$ws.Range("E42").Value = '  +0.26%  '
$ws.Range("E43").Value = '  -0.84%  '
$ws.Range("E44").Value = '  +5.14%  '
$ws.Range("B45").Value = 'NEARProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D45").Value = '''4.22'
$ws.Range("E45").Value = '  +3.33%  '
$ws.Range("B46").Value = 'WEMIXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D46").Value = '''2.56'
$ws.Range("E46").Value = '  +5.79%  '
$ws.Range("D47").Value = '''2.35'
$ws.Range("E47").Value = '  +18.29%  '
$ws.Range("D48").Value = '''16.38'
$ws.Range("E48").Value = '  -1.00%  '
$ws.Range("B49").Value = 'PEPE'
$ws.Range("C49").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D49").Value = '''0.0₃0536'
$ws.Range("E49").Value = '  +21.75%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '''22.23'
$ws.Range("E50").Value = '  +1.64%  '
$ws.Range("D51").Value = '''113.06'
$ws.Range("E51").Value = '  +12.59%  '
